# Scheduled price-refresh update across the Excalibur_Profits sheets.
# Updates currentAveragePrice / currentAveragePriceNQ / currentAveragePriceHQ /
# LevePriceNQ / LevePriceHQ / LeveProfitNQ / LeveProfitHQ (columns H-N) for a
# set of leve rows on each job sheet, reflecting refreshed market-board data.

$wb = $excel.ActiveWorkbook

function Set-Or-Clear {
    param(
        $ws,
        [string]$cellRef,
        $value
    )
    if ($null -eq $value) {
        $ws.Range($cellRef).ClearContents()
    } else {
        $ws.Range($cellRef).Value = $value
    }
}

# ---------------- ALC ----------------
$ws = $wb.Worksheets.Item("ALC")

Set-Or-Clear $ws "H6" 17
Set-Or-Clear $ws "I6" 17
Set-Or-Clear $ws "J6" 0
Set-Or-Clear $ws "K6" 51
Set-Or-Clear $ws "L6" 0
Set-Or-Clear $ws "M6" 61
Set-Or-Clear $ws "N6" $null

Set-Or-Clear $ws "H16" 10822.333
Set-Or-Clear $ws "I16" 4972.5
Set-Or-Clear $ws "J16" 13747.25
Set-Or-Clear $ws "K16" 4972.5
Set-Or-Clear $ws "L16" 13747.25
Set-Or-Clear $ws "M16" -4742.5
Set-Or-Clear $ws "N16" -14207.25

Set-Or-Clear $ws "H98" 1254
Set-Or-Clear $ws "I98" 1015.7037
Set-Or-Clear $ws "K98" 1015.7037
Set-Or-Clear $ws "M98" 482.2963

Set-Or-Clear $ws "H107" 573.6667
Set-Or-Clear $ws "I107" 587.53845
Set-Or-Clear $ws "K107" 587.53845
Set-Or-Clear $ws "M107" 1332.46155

Set-Or-Clear $ws "H111" 8538.125
Set-Or-Clear $ws "I111" 5925
Set-Or-Clear $ws "K111" 17775
Set-Or-Clear $ws "M111" -14708

Set-Or-Clear $ws "H112" 1023.1667
Set-Or-Clear $ws "J112" 1022.6667
Set-Or-Clear $ws "L112" 3068.0001
Set-Or-Clear $ws "N112" -5284.0001

Set-Or-Clear $ws "H122" 1254
Set-Or-Clear $ws "I122" 1015.7037
Set-Or-Clear $ws "K122" 3047.1111
Set-Or-Clear $ws "M122" -597.1111000000001

Set-Or-Clear $ws "H137" 38463076
Set-Or-Clear $ws "I137" 62501212
Set-Or-Clear $ws "J137" 2059.9
Set-Or-Clear $ws "K137" 187503636
Set-Or-Clear $ws "L137" 6179.700000000001
Set-Or-Clear $ws "M137" -187501086
Set-Or-Clear $ws "N137" -11279.7

# ---------------- ARM ----------------
$ws = $wb.Worksheets.Item("ARM")

Set-Or-Clear $ws "H4" 500
Set-Or-Clear $ws "J4" 0
Set-Or-Clear $ws "L4" 0
Set-Or-Clear $ws "N4" $null

Set-Or-Clear $ws "H32" 821.9403
Set-Or-Clear $ws "I32" 790.15625
Set-Or-Clear $ws "K32" 790.15625
Set-Or-Clear $ws "M32" -503.15625

Set-Or-Clear $ws "H74" 1604.1765
Set-Or-Clear $ws "I74" 1018.06665
Set-Or-Clear $ws "K74" 1018.06665
Set-Or-Clear $ws "M74" -144.06665

Set-Or-Clear $ws "H77" 1604.1765
Set-Or-Clear $ws "I77" 1018.06665
Set-Or-Clear $ws "K77" 5090.33325
Set-Or-Clear $ws "M77" -722.3332499999997

Set-Or-Clear $ws "H110" 1614.7858
Set-Or-Clear $ws "I110" 1357.3334
Set-Or-Clear $ws "J110" 1807.875
Set-Or-Clear $ws "K110" 1357.3334
Set-Or-Clear $ws "L110" 1807.875
Set-Or-Clear $ws "M110" 687.6666
Set-Or-Clear $ws "N110" -5897.875

Set-Or-Clear $ws "H132" 3298.6191
Set-Or-Clear $ws "I132" 2749.3235
Set-Or-Clear $ws "J132" 5633.125
Set-Or-Clear $ws "K132" 8247.970499999999
Set-Or-Clear $ws "L132" 16899.375
Set-Or-Clear $ws "M132" -5717.970499999999
Set-Or-Clear $ws "N132" -21959.375

# ---------------- BSM ----------------
$ws = $wb.Worksheets.Item("BSM")

Set-Or-Clear $ws "H94" 1361.2
Set-Or-Clear $ws "I94" 1264.5
Set-Or-Clear $ws "J94" 1554.6
Set-Or-Clear $ws "K94" 1264.5
Set-Or-Clear $ws "L94" 1554.6
Set-Or-Clear $ws "M94" -813.5
Set-Or-Clear $ws "N94" -2456.6

Set-Or-Clear $ws "H107" 3863.1428
Set-Or-Clear $ws "I107" 4467.4
Set-Or-Clear $ws "J107" 2352.5
Set-Or-Clear $ws "K107" 4467.4
Set-Or-Clear $ws "L107" 2352.5
Set-Or-Clear $ws "M107" -2547.4
Set-Or-Clear $ws "N107" -6192.5

Set-Or-Clear $ws "H134" 2033.0944
Set-Or-Clear $ws "I134" 1313
Set-Or-Clear $ws "K134" 3939
Set-Or-Clear $ws "M134" -1404

# ---------------- CRP ----------------
$ws = $wb.Worksheets.Item("CRP")

Set-Or-Clear $ws "H35" 12068.125
Set-Or-Clear $ws "I35" 12068.125
Set-Or-Clear $ws "K35" 12068.125
Set-Or-Clear $ws "M35" -11774.125

Set-Or-Clear $ws "H132" 105266790
Set-Or-Clear $ws "I132" 181819550
Set-Or-Clear $ws "J132" 6762.25
Set-Or-Clear $ws "K132" 545458650
Set-Or-Clear $ws "L132" 20286.75
Set-Or-Clear $ws "M132" -545456120
Set-Or-Clear $ws "N132" -25346.75

Set-Or-Clear $ws "H134" 17066.953
Set-Or-Clear $ws "I134" 21570.562
Set-Or-Clear $ws "K134" 64711.686
Set-Or-Clear $ws "M134" -62176.686

# ---------------- CUL ----------------
$ws = $wb.Worksheets.Item("CUL")

Set-Or-Clear $ws "H68" 19603
Set-Or-Clear $ws "I68" 0
Set-Or-Clear $ws "J68" 19603
Set-Or-Clear $ws "K68" 0
Set-Or-Clear $ws "L68" 58809
Set-Or-Clear $ws "M68" $null
Set-Or-Clear $ws "N68" -60431

Set-Or-Clear $ws "H71" 19603
Set-Or-Clear $ws "I71" 0
Set-Or-Clear $ws "J71" 19603
Set-Or-Clear $ws "K71" 0
Set-Or-Clear $ws "L71" 176427
Set-Or-Clear $ws "M71" $null
Set-Or-Clear $ws "N71" -184539

# ---------------- GSM ----------------
$ws = $wb.Worksheets.Item("GSM")

Set-Or-Clear $ws "H132" 23268674
Set-Or-Clear $ws "I132" 34495324
Set-Or-Clear $ws "J132" 13469.643
Set-Or-Clear $ws "K132" 103485972
Set-Or-Clear $ws "L132" 40408.929
Set-Or-Clear $ws "M132" -103483442
Set-Or-Clear $ws "N132" -45468.929

Set-Or-Clear $ws "H140" 98000
Set-Or-Clear $ws "J140" 98000
Set-Or-Clear $ws "L140" 98000
Set-Or-Clear $ws "N140" -108360

# ---------------- LTW ----------------
$ws = $wb.Worksheets.Item("LTW")

Set-Or-Clear $ws "H55" 400.25
Set-Or-Clear $ws "I55" 300.5
Set-Or-Clear $ws "K55" 300.5
Set-Or-Clear $ws "M55" -127.5

Set-Or-Clear $ws "H61" 1066.1666
Set-Or-Clear $ws "I61" 879.7
Set-Or-Clear $ws "K61" 879.7
Set-Or-Clear $ws "M61" -677.7

Set-Or-Clear $ws "H82" 3250
Set-Or-Clear $ws "I82" 4000
Set-Or-Clear $ws "K82" 4000
Set-Or-Clear $ws "M82" -3639

Set-Or-Clear $ws "H85" 3250
Set-Or-Clear $ws "I85" 4000
Set-Or-Clear $ws "K85" 4000
Set-Or-Clear $ws "M85" -2752

Set-Or-Clear $ws "H113" 1066.1666
Set-Or-Clear $ws "I113" 879.7
Set-Or-Clear $ws "K113" 879.7
Set-Or-Clear $ws "M113" 1290.3

Set-Or-Clear $ws "H122" 4963.9565
Set-Or-Clear $ws "I122" 4535.316
Set-Or-Clear $ws "K122" 13605.948
Set-Or-Clear $ws "M122" -11155.948

Set-Or-Clear $ws "H132" 3635.9092
Set-Or-Clear $ws "I132" 2999.8333
Set-Or-Clear $ws "J132" 4399.2
Set-Or-Clear $ws "K132" 8999.499899999999
Set-Or-Clear $ws "L132" 13197.6
Set-Or-Clear $ws "M132" -6469.499899999999
Set-Or-Clear $ws "N132" -18257.6

# ---------------- WVR ----------------
$ws = $wb.Worksheets.Item("WVR")

Set-Or-Clear $ws "H132" 4655237.5
Set-Or-Clear $ws "I132" 5717211.5
Set-Or-Clear $ws "K132" 17151634.5
Set-Or-Clear $ws "M132" -17149104.5

Set-Or-Clear $ws "H136" 7094909
Set-Or-Clear $ws "I136" 7578415
Set-Or-Clear $ws "K136" 22735245
Set-Or-Clear $ws "M136" -22732695
